$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 47357.5  # H33: 54098.93 -> 47357.5
$ws.Cells.Item(33, 9).Value = 62812.25  # I33: 75341.2 -> 62812.25
$ws.Cells.Item(33, 11).Value = 62812.25  # K33: 75341.2 -> 62812.25
$ws.Cells.Item(33, 13).Value = -62583.25  # M33: -75112.2 -> -62583.25

$ws.Cells.Item(41, 8).Value = 77245.234  # H41: 91285.27 -> 77245.234
$ws.Cells.Item(41, 9).Value = 375  # I41: 600 -> 375
$ws.Cells.Item(41, 10).Value = 91221.63  # J41: 100353.8 -> 91221.63
$ws.Cells.Item(41, 11).Value = 375  # K41: 600 -> 375
$ws.Cells.Item(41, 12).Value = 91221.63  # L41: 100353.8 -> 91221.63
$ws.Cells.Item(41, 13).Value = 65  # M41: -160 -> 65
$ws.Cells.Item(41, 14).Value = -92101.63  # N41: -101233.8 -> -92101.63

$ws.Cells.Item(62, 8).Value = 1081609.8  # H62: 1081689.8 -> 1081609.8
$ws.Cells.Item(62, 10).Value = 126674.75  # J62: 126874.75 -> 126674.75
$ws.Cells.Item(62, 12).Value = 126674.75  # L62: 126874.75 -> 126674.75
$ws.Cells.Item(62, 14).Value = -127922.75  # N62: -128122.75 -> -127922.75

$ws.Cells.Item(65, 8).Value = 1081609.8  # H65: 1081689.8 -> 1081609.8
$ws.Cells.Item(65, 10).Value = 126674.75  # J65: 126874.75 -> 126674.75
$ws.Cells.Item(65, 12).Value = 633373.75  # L65: 634373.75 -> 633373.75
$ws.Cells.Item(65, 14).Value = -639613.75  # N65: -640613.75 -> -639613.75

$ws.Cells.Item(86, 8).Value = 20105030  # H86: 20105078 -> 20105030
$ws.Cells.Item(86, 9).Value = 4560  # I86: 4465.8335 -> 4560
$ws.Cells.Item(86, 10).Value = 40205500  # J86: 50255996 -> 40205500
$ws.Cells.Item(86, 11).Value = 4560  # K86: 4465.8335 -> 4560
$ws.Cells.Item(86, 12).Value = 40205500  # L86: 50255996 -> 40205500
$ws.Cells.Item(86, 13).Value = -3437  # M86: -3342.8335 -> -3437
$ws.Cells.Item(86, 14).Value = -40207746  # N86: -50258242 -> -40207746

$ws.Cells.Item(89, 8).Value = 20105030  # H89: 20105078 -> 20105030
$ws.Cells.Item(89, 9).Value = 4560  # I89: 4465.8335 -> 4560
$ws.Cells.Item(89, 10).Value = 40205500  # J89: 50255996 -> 40205500
$ws.Cells.Item(89, 11).Value = 22800  # K89: 22329.1675 -> 22800
$ws.Cells.Item(89, 12).Value = 201027500  # L89: 251279980 -> 201027500
$ws.Cells.Item(89, 13).Value = -17184  # M89: -16713.1675 -> -17184
$ws.Cells.Item(89, 14).Value = -201038732  # N89: -251291212 -> -201038732

$ws.Cells.Item(98, 8).Value = 391.18182  # H98: 398.85715 -> 391.18182
$ws.Cells.Item(98, 9).Value = 391.18182  # I98: 398.85715 -> 391.18182
$ws.Cells.Item(98, 11).Value = 391.18182  # K98: 398.85715 -> 391.18182
$ws.Cells.Item(98, 13).Value = 1106.81818  # M98: 1099.14285 -> 1106.81818

$ws.Cells.Item(122, 8).Value = 391.18182  # H122: 398.85715 -> 391.18182
$ws.Cells.Item(122, 9).Value = 391.18182  # I122: 398.85715 -> 391.18182
$ws.Cells.Item(122, 11).Value = 1173.54546  # K122: 1196.57145 -> 1173.54546
$ws.Cells.Item(122, 13).Value = 1276.45454  # M122: 1253.42855 -> 1276.45454

$ws.Cells.Item(137, 8).Value = 1741.5294  # H137: 1765.0883 -> 1741.5294
$ws.Cells.Item(137, 10).Value = 2592.4546  # J137: 2665.2727 -> 2592.4546
$ws.Cells.Item(137, 12).Value = 7777.3638  # L137: 7995.8181 -> 7777.3638
$ws.Cells.Item(137, 14).Value = -12877.3638  # N137: -13095.8181 -> -12877.3638

$ws.Cells.Item(138, 8).Value = 1955.6666  # H138: 1911.775 -> 1955.6666
$ws.Cells.Item(138, 9).Value = 1213.826  # I138: 1145.44 -> 1213.826
$ws.Cells.Item(138, 10).Value = 3268.1538  # J138: 3189 -> 3268.1538
$ws.Cells.Item(138, 11).Value = 3641.478  # K138: 3436.32 -> 3641.478
$ws.Cells.Item(138, 12).Value = 9804.4614  # L138: 9567 -> 9804.4614
$ws.Cells.Item(138, 13).Value = 1498.522  # M138: 1703.68 -> 1498.522
$ws.Cells.Item(138, 14).Value = -20084.4614  # N138: -19847 -> -20084.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2553.029  # H32: 2633.1594 -> 2553.029
$ws.Cells.Item(32, 9).Value = 1392.8167  # I32: 1425.5593 -> 1392.8167
$ws.Cells.Item(32, 10).Value = 10287.777  # J32: 9758 -> 10287.777
$ws.Cells.Item(32, 11).Value = 1392.8167  # K32: 1425.5593 -> 1392.8167
$ws.Cells.Item(32, 12).Value = 10287.777  # L32: 9758 -> 10287.777
$ws.Cells.Item(32, 13).Value = -1105.8167  # M32: -1138.5593 -> -1105.8167
$ws.Cells.Item(32, 14).Value = -10861.777  # N32: -10332 -> -10861.777

$ws.Cells.Item(61, 8).Value = 37039452  # H61: 47621724 -> 37039452
$ws.Cells.Item(61, 9).Value = 55556772  # I61: 83334410 -> 55556772
$ws.Cells.Item(61, 11).Value = 55556772  # K61: 83334410 -> 55556772
$ws.Cells.Item(61, 13).Value = -55556560  # M61: -83334198 -> -55556560

$ws.Cells.Item(102, 8).Value = 2682325.8  # H102: 2941910 -> 2682325.8
$ws.Cells.Item(102, 9).Value = 3031575.5  # I102: 3136153 -> 3031575.5
$ws.Cells.Item(102, 10).Value = 62953.5  # J102: 125387 -> 62953.5
$ws.Cells.Item(102, 11).Value = 3031575.5  # K102: 3136153 -> 3031575.5
$ws.Cells.Item(102, 12).Value = 62953.5  # L102: 125387 -> 62953.5
$ws.Cells.Item(102, 13).Value = -3029953.5  # M102: -3134531 -> -3029953.5
$ws.Cells.Item(102, 14).Value = -66197.5  # N102: -128631 -> -66197.5

$ws.Cells.Item(122, 8).Value = 8774836  # H122: 8336201.5 -> 8774836
$ws.Cells.Item(122, 9).Value = 10419368  # I122: 9806591 -> 10419368
$ws.Cells.Item(122, 11).Value = 31258104  # K122: 29419773 -> 31258104
$ws.Cells.Item(122, 13).Value = -31255654  # M122: -29417323 -> -31255654

$ws.Cells.Item(132, 8).Value = 58826080  # H132: 35716268 -> 58826080
$ws.Cells.Item(132, 9).Value = 90911550  # I132: 45456324 -> 90911550
$ws.Cells.Item(132, 11).Value = 272734650  # K132: 136368972 -> 272734650
$ws.Cells.Item(132, 13).Value = -272732120  # M132: -136366442 -> -272732120

$ws.Cells.Item(136, 8).Value = 37039452  # H136: 47621724 -> 37039452
$ws.Cells.Item(136, 9).Value = 55556772  # I136: 83334410 -> 55556772
$ws.Cells.Item(136, 11).Value = 166670316  # K136: 250003230 -> 166670316
$ws.Cells.Item(136, 13).Value = -166667766  # M136: -250000680 -> -166667766

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 7750  # H26: 8250 -> 7750
$ws.Cells.Item(26, 9).Value = 7750  # I26: 8250 -> 7750
$ws.Cells.Item(26, 11).Value = 7750  # K26: 8250 -> 7750
$ws.Cells.Item(26, 13).Value = -7458  # M26: -7958 -> -7458

$ws.Cells.Item(86, 8).Value = 7024.125  # H86: 7570.4287 -> 7024.125
$ws.Cells.Item(86, 9).Value = 7949.8335  # I86: 8899.799999999999 -> 7949.8335
$ws.Cells.Item(86, 11).Value = 7949.8335  # K86: 8899.799999999999 -> 7949.8335
$ws.Cells.Item(86, 13).Value = -6826.8335  # M86: -7776.799999999999 -> -6826.8335

$ws.Cells.Item(89, 8).Value = 7024.125  # H89: 7570.4287 -> 7024.125
$ws.Cells.Item(89, 9).Value = 7949.8335  # I89: 8899.799999999999 -> 7949.8335
$ws.Cells.Item(89, 11).Value = 39749.1675  # K89: 44499 -> 39749.1675
$ws.Cells.Item(89, 13).Value = -34133.1675  # M89: -38883 -> -34133.1675

$ws.Cells.Item(94, 8).Value = 5147.647  # H94: 4916.9443 -> 5147.647
$ws.Cells.Item(94, 9).Value = 6547.364  # I94: 6084.6665 -> 6547.364
$ws.Cells.Item(94, 11).Value = 6547.364  # K94: 6084.6665 -> 6547.364
$ws.Cells.Item(94, 13).Value = -6096.364  # M94: -5633.6665 -> -6096.364

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 170.93103  # H22: 175.82143 -> 170.93103
$ws.Cells.Item(22, 9).Value = 198.71428  # I22: 214.6842 -> 198.71428
$ws.Cells.Item(22, 10).Value = 98  # J22: 93.77778000000001 -> 98
$ws.Cells.Item(22, 11).Value = 198.71428  # K22: 214.6842 -> 198.71428
$ws.Cells.Item(22, 12).Value = 98  # L22: 93.77778000000001 -> 98
$ws.Cells.Item(22, 13).Value = 151.28572  # M22: 135.3158 -> 151.28572
$ws.Cells.Item(22, 14).Value = -798  # N22: -793.77778 -> -798

$ws.Cells.Item(58, 8).Value = 1941.05  # H58: 1909.238 -> 1941.05
$ws.Cells.Item(58, 9).Value = 1840.6428  # I58: 1850.5 -> 1840.6428
$ws.Cells.Item(58, 10).Value = 2175.3333  # J58: 2026.7142 -> 2175.3333
$ws.Cells.Item(58, 11).Value = 1840.6428  # K58: 1850.5 -> 1840.6428
$ws.Cells.Item(58, 12).Value = 2175.3333  # L58: 2026.7142 -> 2175.3333
$ws.Cells.Item(58, 13).Value = -1637.6428  # M58: -1647.5 -> -1637.6428
$ws.Cells.Item(58, 14).Value = -2581.3333  # N58: -2432.7142 -> -2581.3333

$ws.Cells.Item(107, 8).Value = 3570.4  # H107: 3811.5557 -> 3570.4
$ws.Cells.Item(107, 9).Value = 2567.5334  # I107: 2650.9285 -> 2567.5334
$ws.Cells.Item(107, 10).Value = 6579  # J107: 7873.75 -> 6579
$ws.Cells.Item(107, 11).Value = 2567.5334  # K107: 2650.9285 -> 2567.5334
$ws.Cells.Item(107, 12).Value = 6579  # L107: 7873.75 -> 6579
$ws.Cells.Item(107, 13).Value = -647.5333999999998  # M107: -730.9285 -> -647.5333999999998
$ws.Cells.Item(107, 14).Value = -10419  # N107: -11713.75 -> -10419

$ws.Cells.Item(132, 8).Value = 2453.182  # H132: 2471.1365 -> 2453.182
$ws.Cells.Item(132, 10).Value = 2749.5  # J132: 2947 -> 2749.5
$ws.Cells.Item(132, 12).Value = 8248.5  # L132: 8841 -> 8248.5
$ws.Cells.Item(132, 14).Value = -13308.5  # N132: -13901 -> -13308.5

$ws.Cells.Item(136, 8).Value = 1941.05  # H136: 1909.238 -> 1941.05
$ws.Cells.Item(136, 9).Value = 1840.6428  # I136: 1850.5 -> 1840.6428
$ws.Cells.Item(136, 10).Value = 2175.3333  # J136: 2026.7142 -> 2175.3333
$ws.Cells.Item(136, 11).Value = 5521.928400000001  # K136: 5551.5 -> 5521.928400000001
$ws.Cells.Item(136, 12).Value = 6525.999899999999  # L136: 6080.142599999999 -> 6525.999899999999
$ws.Cells.Item(136, 13).Value = -2971.928400000001  # M136: -3001.5 -> -2971.928400000001
$ws.Cells.Item(136, 14).Value = -11625.9999  # N136: -11180.1426 -> -11625.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 6565.4443  # H76: 7471.4165 -> 6565.4443

$ws.Cells.Item(79, 8).Value = 6565.4443  # H79: 7471.4165 -> 6565.4443

$ws.Cells.Item(109, 8).Value = 289.875  # H109: 325.7 -> 289.875
$ws.Cells.Item(109, 9).Value = 295.7143  # I109: 352.85715 -> 295.7143
$ws.Cells.Item(109, 10).Value = 249  # J109: 262.33334 -> 249
$ws.Cells.Item(109, 11).Value = 887.1428999999999  # K109: 1058.57145 -> 887.1428999999999
$ws.Cells.Item(109, 12).Value = 747  # L109: 787.0000200000001 -> 747
$ws.Cells.Item(109, 13).Value = 152.8571000000001  # M109: -18.57144999999991 -> 152.8571000000001
$ws.Cells.Item(109, 14).Value = -2827  # N109: -2867.00002 -> -2827

$ws.Cells.Item(112, 8).Value = 3215.2856  # H112: 4515 -> 3215.2856
$ws.Cells.Item(112, 9).Value = 3412.8333  # I112: 4515 -> 3412.8333
$ws.Cells.Item(112, 10).Value = 2030  # J112: 0 -> 2030
$ws.Cells.Item(112, 11).Value = 10238.4999  # K112: 13545 -> 10238.4999
$ws.Cells.Item(112, 12).Value = 6090  # L112: 0 -> 6090
$ws.Cells.Item(112, 13).Value = -9130.499899999999  # M112: -12437 -> -9130.499899999999
$ws.Cells.Item(112, 14).Value = -8306  # N112: None -> -8306

$ws.Cells.Item(115, 8).Value = 3014.2856  # H115: 3148.1667 -> 3014.2856
$ws.Cells.Item(115, 9).Value = 2211  # I115: 0 -> 2211
$ws.Cells.Item(115, 11).Value = 6633  # K115: 0 -> 6633
$ws.Cells.Item(115, 13).Value = -5458  # M115: None -> -5458

$ws.Cells.Item(118, 8).Value = 450  # H118: 0 -> 450
$ws.Cells.Item(118, 9).Value = 450  # I118: 0 -> 450
$ws.Cells.Item(118, 11).Value = 1350  # K118: 0 -> 1350
$ws.Cells.Item(118, 13).Value = -107  # M118: None -> -107

$ws.Cells.Item(121, 8).Value = 819474.9399999999  # H121: 936452.5600000001 -> 819474.9399999999
$ws.Cells.Item(121, 10).Value = 100821.8  # J121: 125869.375 -> 100821.8
$ws.Cells.Item(121, 12).Value = 302465.4  # L121: 377608.125 -> 302465.4
$ws.Cells.Item(121, 14).Value = -305085.4  # N121: -380228.125 -> -305085.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 0  # H6: 1008 -> 0
$ws.Cells.Item(6, 9).Value = 0  # I6: 1008 -> 0
$ws.Cells.Item(6, 11).Value = 0  # K6: 1008 -> 0

$ws.Cells.Item(16, 8).Value = 0  # H16: 1008 -> 0
$ws.Cells.Item(16, 9).Value = 0  # I16: 1008 -> 0
$ws.Cells.Item(16, 11).Value = 0  # K16: 1008 -> 0

$ws.Cells.Item(70, 8).Value = 8255.286  # H70: 8550.571 -> 8255.286
$ws.Cells.Item(70, 9).Value = 7834.9375  # I70: 8037.3335 -> 7834.9375
$ws.Cells.Item(70, 10).Value = 9600.4  # J70: 9833.666999999999 -> 9600.4
$ws.Cells.Item(70, 11).Value = 7834.9375  # K70: 8037.3335 -> 7834.9375
$ws.Cells.Item(70, 12).Value = 9600.4  # L70: 9833.666999999999 -> 9600.4
$ws.Cells.Item(70, 13).Value = -7564.9375  # M70: -7767.3335 -> -7564.9375
$ws.Cells.Item(70, 14).Value = -10140.4  # N70: -10373.667 -> -10140.4

$ws.Cells.Item(73, 8).Value = 8255.286  # H73: 8550.571 -> 8255.286
$ws.Cells.Item(73, 9).Value = 7834.9375  # I73: 8037.3335 -> 7834.9375
$ws.Cells.Item(73, 10).Value = 9600.4  # J73: 9833.666999999999 -> 9600.4
$ws.Cells.Item(73, 11).Value = 7834.9375  # K73: 8037.3335 -> 7834.9375
$ws.Cells.Item(73, 12).Value = 9600.4  # L73: 9833.666999999999 -> 9600.4
$ws.Cells.Item(73, 13).Value = -6898.9375  # M73: -7101.3335 -> -6898.9375
$ws.Cells.Item(73, 14).Value = -11472.4  # N73: -11705.667 -> -11472.4

$ws.Cells.Item(110, 8).Value = 0  # H110: 44478 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 44478 -> 0
$ws.Cells.Item(110, 14).Value = 0  # N110: -52658 -> 0

$ws.Cells.Item(113, 8).Value = 2227.8667  # H113: 2275.5 -> 2227.8667
$ws.Cells.Item(113, 10).Value = 2800.5  # J113: 2821.5557 -> 2800.5
$ws.Cells.Item(113, 12).Value = 2800.5  # L113: 2821.5557 -> 2800.5
$ws.Cells.Item(113, 14).Value = -7140.5  # N113: -7161.5557 -> -7140.5

$ws.Cells.Item(122, 8).Value = 2661.9333  # H122: 2411.85 -> 2661.9333
$ws.Cells.Item(122, 9).Value = 2335.9167  # I122: 2170.8125 -> 2335.9167
$ws.Cells.Item(122, 10).Value = 3966  # J122: 3376 -> 3966
$ws.Cells.Item(122, 11).Value = 7007.750100000001  # K122: 6512.4375 -> 7007.750100000001
$ws.Cells.Item(122, 12).Value = 11898  # L122: 10128 -> 11898
$ws.Cells.Item(122, 13).Value = -4557.750100000001  # M122: -4062.4375 -> -4557.750100000001
$ws.Cells.Item(122, 14).Value = -16798  # N122: -15028 -> -16798

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2664.5405  # H136: 2612.6155 -> 2664.5405
$ws.Cells.Item(136, 9).Value = 2459.6333  # I136: 2409.1562 -> 2459.6333
$ws.Cells.Item(136, 11).Value = 7378.8999  # K136: 7227.4686 -> 7378.8999
$ws.Cells.Item(136, 13).Value = -4828.8999  # M136: -4677.4686 -> -4828.8999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 38000  # H70: 37000 -> 38000
$ws.Cells.Item(70, 10).Value = 0  # J70: 36000 -> 0
$ws.Cells.Item(70, 12).Value = 0  # L70: 36000 -> 0
$ws.Cells.Item(70, 14).ClearContents()  # N70: -36630 -> (cleared)

$ws.Cells.Item(73, 8).Value = 38000  # H73: 37000 -> 38000
$ws.Cells.Item(73, 10).Value = 0  # J73: 36000 -> 0
$ws.Cells.Item(73, 12).Value = 0  # L73: 36000 -> 0
$ws.Cells.Item(73, 14).ClearContents()  # N73: -38184 -> (cleared)

$ws.Cells.Item(100, 8).Value = 25550  # H100: 9024.5 -> 25550
$ws.Cells.Item(100, 9).Value = 50000  # I100: 17282.334 -> 50000
$ws.Cells.Item(100, 10).Value = 1100  # J100: 766.6667 -> 1100
$ws.Cells.Item(100, 11).Value = 100000  # K100: 34564.668 -> 100000
$ws.Cells.Item(100, 12).Value = 2200  # L100: 1533.3334 -> 2200
$ws.Cells.Item(100, 13).Value = -99459  # M100: -34023.668 -> -99459
$ws.Cells.Item(100, 14).Value = -3282  # N100: -2615.3334 -> -3282

$ws.Cells.Item(132, 8).Value = 4238.607  # H132: 3311.8975 -> 4238.607
$ws.Cells.Item(132, 9).Value = 4780.706  # I132: 3346.5186 -> 4780.706
$ws.Cells.Item(132, 10).Value = 3400.818  # J132: 3234 -> 3400.818
$ws.Cells.Item(132, 11).Value = 14342.118  # K132: 10039.5558 -> 14342.118
$ws.Cells.Item(132, 12).Value = 10202.454  # L132: 9702 -> 10202.454
$ws.Cells.Item(132, 13).Value = -11812.118  # M132: -7509.5558 -> -11812.118
$ws.Cells.Item(132, 14).Value = -15262.454  # N132: -14762 -> -15262.454

$ws.Cells.Item(136, 8).Value = 4765.722  # H136: 4672.7896 -> 4765.722
$ws.Cells.Item(136, 9).Value = 2700.3333  # I136: 2743.1428 -> 2700.3333
$ws.Cells.Item(136, 11).Value = 8100.999899999999  # K136: 8229.428400000001 -> 8100.999899999999
$ws.Cells.Item(136, 13).Value = -5550.999899999999  # M136: -5679.428400000001 -> -5550.999899999999

Write-Output "Applied $(45) row updates across 8 sheets"
